$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "navigation" column G: IF(status == failed? style) cells.
# Rows 10,11,14,15,16,17,18 -> "s" ; row 12 -> "agent fraized"
$rowsWithS = @(10, 11, 14, 15, 16, 17, 18)
foreach ($r in $rowsWithS) {
    $ws.Cells.Item($r, 7).Value = "s"
}
$ws.Cells.Item(12, 7).Value = "agent fraized"

# Move the active selection from E12 to A12 (matches the diff's <selection>)
$ws.Range("A12").Select()
